# This edit rotates the data of rows 7, 8 and 9 on the active sheet:
#   new row 7 <- old row 9
#   new row 8 <- old row 7
#   new row 9 <- old row 8
# (Row positions stay put; only the row *contents* move.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are populated somewhere in rows 7-9 (per the sheet data).
# Y and AA (the date-stamp columns) are deliberately left out: they already
# hold the identical literal text "2026-02-03" in all three rows, so
# rotating them is a content no-op, and round-tripping date-shaped text
# through Value2 would otherwise coerce it into a numeric date serial
# (and pick up a date number format) -- a side effect the diff does not
# call for.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W", `
          "X","Z","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP", `
          "AQ","AR","AS","AT","AU","AV","AW","AX","AY")

# Snapshot the current ("before") values of rows 7, 8 and 9 for every
# relevant column before we start overwriting anything (Value2 avoids the
# Value-getter quirk and round-trips numbers/strings/booleans/nulls cleanly).
$snapshot = @{}
foreach ($row in 7..9) {
    foreach ($col in $cols) {
        $snapshot["$col$row"] = $ws.Range("$col$row").Value2
    }
}

# Mapping describing where each destination row's data comes from.
$sourceForDestination = @{ 7 = 9; 8 = 7; 9 = 8 }

foreach ($destRow in 7..9) {
    $srcRow = $sourceForDestination[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $snapshot["$col$srcRow"]
    }
}
